# Weekly fruit/vegetable price update: a new observation (week) is
# inserted as row 69 ("Fecha" 45093, i.e. 2023-06-16) in the
# "Bruselas (repollito)" price history sheet; every following row shifts
# down by one (old row 69 -> 70, ..., old row 82 -> 83).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push rows 69:82 down to 70:83, leaving a blank row 69 for the new record.
$ws.Rows.Item(69).Insert()

# Populate the newly inserted row 69 with the new weekly observation.
$ws.Range("A69").Value = 6
$ws.Range("B69").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C69").Value = "Metropolitana"
$ws.Range("D69").Value = 45093
$ws.Range("E69").Value = 13
$ws.Range("F69").Value = 100112035
$ws.Range("G69").Value = "Bruselas (repollito)"
$ws.Range("H69").Value = "Sin especificar"
$ws.Range("I69").Value = "Primera"
$ws.Range("J69").Value = 640
$ws.Range("K69").Value = 17000
$ws.Range("L69").Value = 19000
$ws.Range("M69").Value = 17844
$ws.Range("N69").Value = "$/malla 15 kilos"
$ws.Range("O69").Value = "Provincia de Quillota"
$ws.Range("P69").Value = 1190
$ws.Range("Q69").Value = 15
$ws.Range("R69").Value = "Hortaliza"
